$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the bold / centered / thin-boxed format on B1 first ...
$c1 = $ws.Range("B1")
$c1.Font.Bold = $true
$c1.HorizontalAlignment = -4108   # xlCenter
$c1.VerticalAlignment = -4160     # xlTop
$c1.Borders.LineStyle = 1         # xlContinuous
$c1.Borders.Weight = 2            # xlThin

# ... then clone the exact same style onto A2 via copy/paste-special so no
# extra intermediate cell style records get allocated.
$c1.Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

Write-Output "done"
